# Third group meeting update: add "Part 2" handover info and new meeting
# details table (rows 17, 19-23) to the meeting diary sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: Part 2 group members -----------------------------------------
# A17 mirrors the bold header style used by B6/C6/D6 ("Time start" etc.)
$ws.Range("B6").Copy($ws.Range("A17"))
$ws.Range("A17").Value = "Part 2 Group Members:"
$ws.Range("B17").Value = "Minghao Zeng"
$ws.Range("C17").Value = "Tashya Sathyajit"
$ws.Range("D17").Value = "Disha Rathod"

# --- Row 19: "Meeting Details:" section heading ----------------------------
$ws.Range("A5").Copy($ws.Range("A19"))
$ws.Range("A19").Value = "Meeting Details:"

# --- Row 20: table column headers ------------------------------------------
$ws.Range("A6").Copy($ws.Range("A20"))
$ws.Range("A20").Value = "Date"
$ws.Range("B6").Copy($ws.Range("B20"))
$ws.Range("B20").Value = "Time start"
$ws.Range("C6").Copy($ws.Range("C20"))
$ws.Range("C20").Value = "Time end"
$ws.Range("D6").Copy($ws.Range("D20"))
$ws.Range("D20").Value = "Members present"
$ws.Range("E6").Copy($ws.Range("E20"))
$ws.Range("E20").Value = "Discussions"

# --- Row 21: meeting on 12/10/2023 ------------------------------------------
$ws.Range("A7").Copy($ws.Range("A21"))
$ws.Range("A21").Value = 45211
$ws.Range("B7").Copy($ws.Range("B21"))
$ws.Range("B21").Value = 0.6875
$ws.Range("C7").Copy($ws.Range("C21"))
$ws.Range("C21").Value = 0.72916666666666663
$ws.Range("D7").Copy($ws.Range("D21"))
$ws.Range("D21").Value = "All"
$ws.Range("E8").Copy($ws.Range("E21"))
$ws.Range("E21").Value = "Completing handover process from part 1 group member Arindom. "

# --- Row 22: meeting on 15/10/2023 (planning discussion) -------------------
$ws.Range("A7").Copy($ws.Range("A22"))
$ws.Range("A22").Value = 45214
$ws.Range("B7").Copy($ws.Range("B22"))
$ws.Range("B22").Value = 0.70833333333333337
$ws.Range("C7").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 0.72916666666666663
$ws.Range("D7").Copy($ws.Range("D22"))
$ws.Range("D22").Value = "All"
$ws.Range("E8").Copy($ws.Range("E22"))
$ws.Range("E22").Value = "Discussion on how to approach the assignment. All members agreed to complete the whole assignment by 22/10/2023 and compare/merge the results afterward."

# --- Row 23: meeting on 15/10/2023 (repo setup) -----------------------------
$ws.Range("A7").Copy($ws.Range("A23"))
$ws.Range("A23").Value = 45214
$ws.Range("B7").Copy($ws.Range("B23"))
$ws.Range("B23").Value = 0.89583333333333337
$ws.Range("C7").Copy($ws.Range("C23"))
$ws.Range("C23").Value = 0.9375
$ws.Range("D7").Copy($ws.Range("D23"))
$ws.Range("D23").Value = "All"
$ws.Range("E8").Copy($ws.Range("E23"))
$ws.Range("E23").Value = "New files from template repo added to assignment repo"

# --- Update selection to mirror the saved view (E24, nothing selected yet) -
[void]$ws.Range("E24").Select()
